$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$items = @(
    @{Row=2; Name='CABEÇA HOLLOW 6" COM PORTINHOLA CHAVETA'; Qty=4},
    @{Row=3; Name='CABEÇA HOLLOW 7,5" PARAFUSO COM PORTINHOLA'; Qty=2},
    @{Row=4; Name='CABEÇA HOLLOW 6" CHAVETA'; Qty=1},
    @{Row=5; Name='CABEÇA HOLLOW 7,5" PARFUSO'; Qty=2},
    @{Row=6; Name='HASTE HELICOIDAL 4" 1,5M'; Qty=4},
    @{Row=7; Name='HASTE HEICOIDAL 4" 1M'; Qty=4},
    @{Row=8; Name='HASTE HELICOIDAL 6" 1,5M'; Qty=9},
    @{Row=9; Name='HASTE HELICOIDAL 6" 1M'; Qty=4},
    @{Row=10; Name='PONTEIRA HELICOIDAL 4" 1M'; Qty=5},
    @{Row=11; Name='PONTEIRA HELICOIDAL 6" 1,5M'; Qty=5},
    @{Row=12; Name='HASTE MACIÇA'; Qty=0},
    @{Row=13; Name='HASTE INJEÇÃO 1,2M'; Qty=90},
    @{Row=14; Name='HASTE INJEÇÃO 60CM'; Qty=80},
    @{Row=15; Name='HASTE INOX 1M'; Qty=29},
    @{Row=16; Name='PONTEIRA FIXA INJEÇÃO'; Qty=78},
    @{Row=17; Name='PONTEIRA ARTICULADA'; Qty=26},
    @{Row=18; Name='PONTEIRA TRICONICA 4"'; Qty=1},
    @{Row=19; Name='PONTEIRA TRICONICA 6"'; Qty=1},
    @{Row=20; Name='BATEDOR HASTE MACIÇA'; Qty=0},
    @{Row=21; Name='BATEDOR'; Qty=22},
    @{Row=22; Name='SACADOR GDU'; Qty=1},
    @{Row=23; Name='PESCADOR'; Qty=1},
    @{Row=24; Name='MINI SACADOR'; Qty=18},
    @{Row=25; Name='SACADOR DE HASTE '; Qty=3},
    @{Row=26; Name='EXTRATOR DE HASTE'; Qty=3},
    @{Row=27; Name='GARFO HOLLOW'; Qty=1},
    @{Row=28; Name='GARFO HELICOIDAL'; Qty=3},
    @{Row=29; Name='SACADOR DE LINER'; Qty=4},
    @{Row=30; Name='AMOSTRADOR PROBE 1,2M'; Qty=15},
    @{Row=31; Name='AMOSTRADOR PROBE 60CM'; Qty=8},
    @{Row=32; Name='MINI PROBE'; Qty=1},
    @{Row=33; Name='MINI PROBE BI PARTIDO'; Qty=1},
    @{Row=34; Name='ADAPTADOR PROBE MACHO'; Qty=3},
    @{Row=35; Name='ADAPTADOR PROBE FEMEA'; Qty=3},
    @{Row=36; Name='PROBE BI PARTIDO '; Qty=4},
    @{Row=37; Name='TE DE INJEÇÃO'; Qty=4},
    @{Row=38; Name='T INOX COM SAIDA MANOMETRO'; Qty=11},
    @{Row=39; Name='CORTADOR DE LINER'; Qty=2},
    @{Row=40; Name='SUPORTE DE CORTAR LINER'; Qty=3},
    @{Row=41; Name='BAILER DE INOX 1.1/2"'; Qty=1},
    @{Row=42; Name='BAILER D EINOX 3.1/4" '; Qty=2},
    @{Row=43; Name='HOLLOW 6" 1M CHAVETA'; Qty=18},
    @{Row=44; Name='HOLLOW 7,5" 1,5M COM PARAFUSO'; Qty=7},
    @{Row=45; Name='PONTEIRA FIXA INOX'; Qty=0},
    @{Row=46; Name='BOMBA PNEUMATICA'; Qty=0},
    @{Row=47; Name='PONTEIRA HELICOIDAL 6" 1M'; Qty=3}
)

foreach ($item in $items) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Name
    $ws.Cells.Item($item.Row, 2).Value = $item.Qty
}

# Row 48 no longer exists in the final layout (data now spans A1:B47);
# delete the entire row to shrink the used range accordingly.
$ws.Rows.Item(48).Delete()
